$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $rng = $ws.Range($ref)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue $ws "D2" "61.324.82"
Set-TextValue $ws "E2" "  +6.69%  "
Set-TextValue $ws "D3" "2.681.28"
Set-TextValue $ws "E3" "  +10.13%  "
Set-TextValue $ws "D4" "0.995"
Set-TextValue $ws "E4" "  -0.14%  "
Set-TextValue $ws "D5" "513.71"
Set-TextValue $ws "E5" "  +4.83%  "
Set-TextValue $ws "D6" "159.28"
Set-TextValue $ws "E6" "  +2.65%  "
Set-TextValue $ws "D7" "0.616"
Set-TextValue $ws "E7" "  -0.10%  "
Set-TextValue $ws "D8" "0.994"
Set-TextValue $ws "E8" "  -0.16%  "
Set-TextValue $ws "D9" "2.676.53"
Set-TextValue $ws "E9" "  +9.65%  "
Set-TextValue $ws "E10" "  +8.76%  "
Set-TextValue $ws "E11" "  +6.05%  "
Set-TextValue $ws "D12" "0.351"
Set-TextValue $ws "E12" "  +4.09%  "
Set-TextValue $ws "E13" "  +0.82%  "
Set-TextValue $ws "D14" "3.096.79"
Set-TextValue $ws "E14" "  +8.68%  "
Set-TextValue $ws "D15" "60.961.93"
Set-TextValue $ws "E15" "  +6.59%  "
Set-TextValue $ws "D16" "22.34"
Set-TextValue $ws "E16" "  +7.29%  "
Set-TextValue $ws "D17" "0.0000141"
Set-TextValue $ws "E17" "  +5.08%  "
Set-TextValue $ws "D18" "2.663.74"
Set-TextValue $ws "E18" "  +9.03%  "
Set-TextValue $ws "E19" "  +1.25%  "
Set-TextValue $ws "D20" "350.65"
Set-TextValue $ws "E20" "  +6.92%  "
Set-TextValue $ws "E21" "  +5.95%  "
Set-TextValue $ws "D22" "6.20"
Set-TextValue $ws "E22" "  +4.74%  "
Set-TextValue $ws "D23" "0.997"
Set-TextValue $ws "E23" "  -0.15%  "
Set-TextValue $ws "D24" "60.67"
Set-TextValue $ws "E24" "  +4.15%  "
Set-TextValue $ws "D25" "0.426"
Set-TextValue $ws "E25" "  +3.63%  "
Set-TextValue $ws "D26" "2.758.42"
Set-TextValue $ws "E26" "  +8.55%  "
Set-TextValue $ws "E27" "  +4.34%  "
Set-TextValue $ws "D28" "0.992"
Set-TextValue $ws "E28" "  -0.56%  "
Set-TextValue $ws "D29" "0.0₃0871"
Set-TextValue $ws "E29" "  +10.05%  "
Set-TextValue $ws "D30" "7.58"
Set-TextValue $ws "E30" "  +3.42%  "
Set-TextValue $ws "D31" "0.998"
Set-TextValue $ws "E31" "  -0.08%  "
Set-TextValue $ws "E32" "  +5.21%  "
Set-TextValue $ws "D33" "157.51"
Set-TextValue $ws "E33" "  +5.31%  "
Set-TextValue $ws "E34" "  +4.01%  "
Set-TextValue $ws "D35" "5.74"
Set-TextValue $ws "E35" "  +7.82%  "
Set-TextValue $ws "D36" "4.11"
Set-TextValue $ws "E36" "  +10.51%  "
Set-TextValue $ws "E37" "  +6.56%  "
Set-TextValue $ws "B38" "Fetch.AI"
Set-TextValue $ws "C38" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws "D38" "0.884"
Set-TextValue $ws "E38" "  +2.89%  "
Set-TextValue $ws "B39" "Stacks"
Set-TextValue $ws "C39" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws "D39" "1.54"
Set-TextValue $ws "E39" "  +11.35%  "
Set-TextValue $ws "D40" "310.95"
Set-TextValue $ws "E40" "  +16.43%  "
Set-TextValue $ws "D41" "3.80"
Set-TextValue $ws "E41" "  +7.80%  "
Set-TextValue $ws "D42" "0.838"
Set-TextValue $ws "E42" "  +30.35%  "
Set-TextValue $ws "D43" "35.67"
Set-TextValue $ws "E43" "  +4.14%  "
Set-TextValue $ws "D44" "0.648"
Set-TextValue $ws "E44" "  +8.54%  "
Set-TextValue $ws "D45" "0.0579"
Set-TextValue $ws "E45" "  +7.93%  "
Set-TextValue $ws "E46" "  -1.03%  "
Set-TextValue $ws "D47" "20.26"
Set-TextValue $ws "E47" "  +15.70%  "
Set-TextValue $ws "B48" "FirstDigitalUSD"
Set-TextValue $ws "C48" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws "D48" "1.00"
Set-TextValue $ws "E48" "  +0.27%  "
Set-TextValue $ws "B49" "RenderToken"
Set-TextValue $ws "C49" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws "D49" "5.04"
Set-TextValue $ws "E49" "  +7.98%  "
Set-TextValue $ws "D50" "0.0237"
Set-TextValue $ws "E50" "  +3.59%  "
Set-TextValue $ws "D51" "2.032.68"
Set-TextValue $ws "E51" "  +9.23%  "

Write-Output "Applied 98 cell updates"
